$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '35.166.28'
$ws.Cells.Item(2, 5).Value = '  -0.27%  '

$ws.Cells.Item(3, 4).Value = '1.898.65'
$ws.Cells.Item(3, 5).Value = '  -0.25%  '

$ws.Cells.Item(4, 5).Value = '  -0.10%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '252.28'
$ws.Cells.Item(5, 5).Value = '  +2.47%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.693'
$ws.Cells.Item(6, 5).Value = '  -0.07%  '

$ws.Cells.Item(7, 5).Value = '  -0.08%  '

$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '40.85'
$ws.Cells.Item(8, 5).Value = '  -3.16%  '

$ws.Cells.Item(9, 5).Value = '  +2.70%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '52.91'
$ws.Cells.Item(10, 5).Value = '  -0.62%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0754'
$ws.Cells.Item(11, 5).Value = '  +3.56%  '

$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.0984'
$ws.Cells.Item(12, 5).Value = '  -1.20%  '

$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '13.01'
$ws.Cells.Item(13, 5).Value = '  +5.65%  '

$ws.Cells.Item(14, 4).Value = '2.176.14'
$ws.Cells.Item(14, 5).Value = '  -0.14%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.736'
$ws.Cells.Item(15, 5).Value = '  +3.87%  '

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '4.94'
$ws.Cells.Item(16, 5).Value = '  +1.83%  '

$ws.Cells.Item(17, 4).Value = '1.912.85'
$ws.Cells.Item(17, 5).Value = '  +0.55%  '

$ws.Cells.Item(18, 4).Value = '35.173.55'
$ws.Cells.Item(18, 5).Value = '  -0.22%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '73.62'
$ws.Cells.Item(19, 5).Value = '  +1.56%  '

$ws.Cells.Item(20, 4).Value = '0.0₃0833'
$ws.Cells.Item(20, 5).Value = '  +1.11%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '242.02'
$ws.Cells.Item(21, 5).Value = '  +0.24%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '12.91'
$ws.Cells.Item(22, 5).Value = '  +2.30%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '5.03'
$ws.Cells.Item(23, 5).Value = '  +3.76%  '

$ws.Cells.Item(24, 5).Value = '  -0.16%  '

$ws.Cells.Item(25, 5).Value = '  +3.90%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '2.31'
$ws.Cells.Item(26, 5).Value = '  -0.71%  '

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '166.61'
$ws.Cells.Item(27, 5).Value = '  -1.98%  '

$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '8.57'
$ws.Cells.Item(28, 5).Value = '  -0.06%  '

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '18.45'
$ws.Cells.Item(29, 5).Value = '  +0.30%  '

$ws.Cells.Item(30, 5).Value = '  -1.58%  '

$ws.Cells.Item(31, 4).Value = '4.127.93'
$ws.Cells.Item(31, 5).Value = '  -0.56%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '2.11'
$ws.Cells.Item(32, 5).Value = '  +17.96%  '

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.0606'
$ws.Cells.Item(33, 5).Value = '  +5.82%  '

$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '4.31'
$ws.Cells.Item(34, 5).Value = '  +2.73%  '

$ws.Cells.Item(35, 5).Value = '  +17.63%  '

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '4.18'
$ws.Cells.Item(36, 5).Value = '  +1.23%  '

$ws.Cells.Item(37, 5).Value = '  +0.00%  '

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.853'
$ws.Cells.Item(38, 5).Value = '  -12.93%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '1.99'
$ws.Cells.Item(39, 5).Value = '  -2.64%  '

$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '99.95'
$ws.Cells.Item(40, 5).Value = '  +10.35%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '17.06'
$ws.Cells.Item(41, 5).Value = '  +4.67%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.0213'
$ws.Cells.Item(42, 5).Value = '  +1.38%  '

$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '1.10'
$ws.Cells.Item(43, 5).Value = '  -1.02%  '

$ws.Cells.Item(44, 5).Value = '  -4.94%  '

$ws.Cells.Item(45, 4).Value = '1.324.92'
$ws.Cells.Item(45, 5).Value = '  -1.45%  '

$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '2.42'
$ws.Cells.Item(46, 5).Value = '  -0.70%  '

$ws.Cells.Item(47, 5).Value = '  +0.61%  '

$ws.Cells.Item(48, 5).Value = '  -1.57%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '6.59'
$ws.Cells.Item(49, 5).Value = '  +0.15%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '11.93'
$ws.Cells.Item(50, 5).Value = '  -6.27%  '

$ws.Cells.Item(51, 2).Value = 'Cronos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.0740'
$ws.Cells.Item(51, 5).Value = '  +4.51%  '
